$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Address" sheet: billing address (row 2) gets a new street + new phone
# numbers, and the shipping address (row 5) is switched from Steve Jobs /
# Apple to William Gates / MicroSoft.
# ---------------------------------------------------------------------------
$wsAddr = $wb.Worksheets.Item("Address")

# Billing address (row 2)
$wsAddr.Cells.Item(2, 4).Value = "Apple Avenue"
$wsAddr.Cells.Item(2, 6).Value = 789
$wsAddr.Cells.Item(2, 7).Value = 55511234

# Shipping address (row 5)
$wsAddr.Cells.Item(5, 2).Value = "William"
$wsAddr.Cells.Item(5, 3).Value = "Gates"
$wsAddr.Cells.Item(5, 4).Value = "Microsoft Street"
$wsAddr.Cells.Item(5, 6).Value = "MicroSoft"

# Street Address column now needs to fit the longer text.
$wsAddr.Columns.Item(4).ColumnWidth = 14.14

# ---------------------------------------------------------------------------
# "Orders" sheet: new invoice run (Order #512695). One line item (Ipoh
# Coffee) came back out of stock, and two more line items were appended by
# the order-processing queue.
# ---------------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Orders")
$currencyFmt = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

# Row 3 (Ipoh Coffee) failed this run - out of stock, price cells cleared.
$wsOrders.Cells.Item(3, 3).Value = "Unsuccessful"
$wsOrders.Cells.Item(3, 4).Value = "Item is out of stock (zero value). Order was not placed"
$wsOrders.Cells.Item(3, 5).ClearContents()
$wsOrders.Cells.Item(3, 6).ClearContents()

# Remaining successful rows now reference the new order number.
$wsOrders.Cells.Item(2, 4).Value = "Order #512695"
$wsOrders.Cells.Item(4, 4).Value = "Order #512695"
$wsOrders.Cells.Item(5, 4).Value = "Order #512695"
$wsOrders.Cells.Item(6, 4).Value = "Order #512695"
$wsOrders.Cells.Item(7, 4).Value = "Order #512695"
$wsOrders.Cells.Item(8, 4).Value = "Order #512695"

# Two new line items appended to the queue (rows 10-11), plus two trailing
# blank (but currency-formatted) rows left by the process (12-13).
$wsOrders.Cells.Item(10, 3).Value = "Successful"
$wsOrders.Cells.Item(10, 4).Value = "Order #512695"
$wsOrders.Cells.Item(10, 5).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(10, 5).Value = 7
$wsOrders.Cells.Item(10, 6).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(10, 6).Value = 14

$wsOrders.Cells.Item(11, 3).Value = "Successful"
$wsOrders.Cells.Item(11, 4).Value = "Order #512695"
$wsOrders.Cells.Item(11, 5).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(11, 5).Value = 7
$wsOrders.Cells.Item(11, 6).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(11, 6).Value = 164

$wsOrders.Cells.Item(12, 5).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(12, 6).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(13, 5).NumberFormat = $currencyFmt
$wsOrders.Cells.Item(13, 6).NumberFormat = $currencyFmt

# New Product column is wider now that the sheet has grown.
$wsOrders.Columns.Item(1).ColumnWidth = 24.8

# Workbook was saved with the whole sheet selected (Ctrl+A) instead of the
# old D9 cell selection.
[void]$wsOrders.Cells.Select()
